$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format first so numeric-looking strings (e.g. "591.90")
# are stored as text, matching the original inlineStr cells (not auto-converted to numbers).
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '63.852.44'
$ws.Range('E2').Value = '  +5.65%  '
$ws.Range('D3').Value = '2.735.35'
$ws.Range('E3').Value = '  +4.64%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '591.90'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = '151.71'
$ws.Range('E6').Value = '  +5.78%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  +2.42%  '
$ws.Range('D9').Value = '2.768.81'
$ws.Range('E9').Value = '  +5.57%  '
$ws.Range('D10').Value = '6.76'
$ws.Range('E10').Value = '  +3.88%  '
$ws.Range('E11').Value = '  +8.06%  '
$ws.Range('E12').Value = '  +4.21%  '
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').Value = '3.223.96'
$ws.Range('E14').Value = '  +4.74%  '
$ws.Range('D15').Value = '26.58'
$ws.Range('E15').Value = '  +7.62%  '
$ws.Range('D16').Value = '63.743.40'
$ws.Range('E16').Value = '  +5.49%  '
$ws.Range('E17').Value = '  +8.35%  '
$ws.Range('D18').Value = '2.757.86'
$ws.Range('E18').Value = '  +5.35%  '
$ws.Range('D19').Value = '12.05'
$ws.Range('E19').Value = '  +5.99%  '
$ws.Range('D20').Value = '4.89'
$ws.Range('E20').Value = '  +5.56%  '
$ws.Range('D21').Value = '366.39'
$ws.Range('E21').Value = '  +5.61%  '
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '0.536'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').Value = '65.94'
$ws.Range('E25').Value = '  +3.86%  '
$ws.Range('E26').Value = '  +4.72%  '
$ws.Range('D27').Value = '8.67'
$ws.Range('E27').Value = '  +8.53%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').Value = '0.0₃0879'
$ws.Range('E29').Value = '  +10.11%  '
$ws.Range('E30').Value = '  +7.02%  '
$ws.Range('D31').Value = '7.13'
$ws.Range('E31').Value = '  +10.59%  '
$ws.Range('D32').Value = '170.02'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '0.996'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.19'
$ws.Range('E34').Value = '  +18.20%  '
$ws.Range('D35').Value = '20.61'
$ws.Range('E35').Value = '  +5.69%  '
$ws.Range('D36').Value = '4.77'
$ws.Range('E36').Value = '  +11.27%  '
$ws.Range('E37').Value = '  +10.05%  '
$ws.Range('D38').Value = '1.81'
$ws.Range('E38').Value = '  +10.14%  '
$ws.Range('E39').Value = '  +20.66%  '
$ws.Range('D40').Value = '352.72'
$ws.Range('E40').Value = '  +10.48%  '
$ws.Range('D41').Value = '4.25'
$ws.Range('E41').Value = '  +8.92%  '
$ws.Range('D42').Value = '39.42'
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('D43').Value = '5.67'
$ws.Range('E43').Value = '  +13.11%  '
$ws.Range('D44').Value = '22.26'
$ws.Range('E44').Value = '  +11.53%  '
$ws.Range('D45').Value = '145.16'
$ws.Range('E45').Value = '  +6.39%  '
$ws.Range('D46').Value = '22.04'
$ws.Range('E46').Value = '  +9.88%  '
$ws.Range('D47').Value = '0.0594'
$ws.Range('E47').Value = '  +7.94%  '
$ws.Range('E48').Value = '  +5.70%  '
$ws.Range('E49').Value = '  +7.41%  '
$ws.Range('E50').Value = '  +2.38%  '
$ws.Range('D51').Value = '2.174.27'
$ws.Range('E51').Value = '  +7.20%  '

# Restore default styling (no explicit style index) now that text is committed,
# matching the original cells which carried no "s" attribute.
$priceRange.Style = "Normal"
